$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.943.95"
$ws.Range("D3").Value = "1.825.37"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9967"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.82"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6144"
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9968"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07413"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2907"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.92"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07625"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.826.80"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.970"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6694"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.44"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009145"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.874"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "28.922.95"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "2.075.80"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.93"
$ws.Range("E20").Value = "  +6.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.61"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9970"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.173"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9977"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.98"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1402"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.454"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.78"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.492"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05563"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.121"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.090"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.196"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.830"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7368"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.135"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.646"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.763"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01775"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "1.206.21"
$ws.Range("E40").Value = "  -2.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.402"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8911"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9949"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.96"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.976.38"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.07"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5071"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4039"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.110"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05801"
$ws.Range("E51").Value = "  +0.27%  "
